$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove Hannah Irwin's result row (row 6), shifting subsequent rows up
$ws.Rows.Item(6).Delete()

# Update the selection to match the post-edit state
$ws.Range("A19").Select()
